$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new rows 20-28 with Arabic ("ara") language reason_list entries ---
# (mirrors the existing "eng"/"fra" blocks already present in rows 2-19)
$ws.Range("A20").Value = "ara"
$ws.Range("B20").Value = "APM"
$ws.Range("C20").Value = "عدم تطابق الصورة مع العمر"
$ws.Range("D20").Value = "عدم التطابق بين العمر والصورة"
$ws.Range("E20").Value = "CLR"
$ws.Range("F20").Value = "'TRUE"

$ws.Range("A21").Value = "ara"
$ws.Range("B21").Value = "GPM"
$ws.Range("C21").Value = "عدم تطابق الصورة بين الجنسين"
$ws.Range("D21").Value = "عدم تطابق الصورة بين الجنسين والصورة"
$ws.Range("E21").Value = "CLR"
$ws.Range("F21").Value = "'TRUE"

$ws.Range("A22").Value = "ara"
$ws.Range("B22").Value = "IAD"
$ws.Range("C22").Value = "عنوان خاطئ"
$ws.Range("D22").Value = "تم العثور على عنوان غير صالح"
$ws.Range("E22").Value = "CLR"
$ws.Range("F22").Value = "'TRUE"

$ws.Range("A23").Value = "ara"
$ws.Range("B23").Value = "DPG"
$ws.Range("C23").Value = "تسجيل مكرر"
$ws.Range("D23").Value = "تم العثور على تسجيل مكرر"
$ws.Range("E23").Value = "CLR"
$ws.Range("F23").Value = "'TRUE"

$ws.Range("A24").Value = "ara"
$ws.Range("B24").Value = "OTH"
$ws.Range("C24").Value = "آحرون"
$ws.Range("D24").Value = "آحرون"
$ws.Range("E24").Value = "CLR"
$ws.Range("F24").Value = "'TRUE"

$ws.Range("A25").Value = "ara"
$ws.Range("B25").Value = "ADM"
$ws.Range("C25").Value = "كل التفاصيل متطابقة"
$ws.Range("D25").Value = "كل التفاصيل متطابقة"
$ws.Range("E25").Value = "MNA"
$ws.Range("F25").Value = "'TRUE"

$ws.Range("A26").Value = "ara"
$ws.Range("B26").Value = "ADD"
$ws.Range("C26").Value = "جميع التفاصيل الديموغرافية متطابقة"
$ws.Range("D26").Value = "جميع التفاصيل الديموغرافية متطابقة"
$ws.Range("E26").Value = "MNA"
$ws.Range("F26").Value = "'TRUE"

$ws.Range("A27").Value = "ara"
$ws.Range("B27").Value = "OPM"
$ws.Range("C27").Value = "فقط الصورة هي المطابقة"
$ws.Range("D27").Value = "فقط الصورة هي المطابقة"
$ws.Range("E27").Value = "MNA"
$ws.Range("F27").Value = "'TRUE"

$ws.Range("A28").Value = "ara"
$ws.Range("B28").Value = "SDM"
$ws.Range("C28").Value = "بعض التفاصيل الديموغرافية متطابقة"
$ws.Range("D28").Value = "بعض التفاصيل الديموغرافية متطابقة"
$ws.Range("E28").Value = "MNA"
$ws.Range("F28").Value = "'TRUE"

# --- Copy the existing header-row cell formatting (style index 4 / 1) onto the new rows ---
# columns A, B, E use the plain data style; column F uses the "text" number format style
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A20:B20").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:B21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:B22").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:B23").PasteSpecial(-4122) | Out-Null
$ws.Range("A24:B24").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:B25").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:B26").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:B27").PasteSpecial(-4122) | Out-Null
$ws.Range("A28:B28").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F24").PasteSpecial(-4122) | Out-Null
$ws.Range("F25").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null

# --- Apply a new wrap-text / left-aligned style to the descr columns (C, D) of the new rows ---
$ws.Range("C20").HorizontalAlignment = -4131
$ws.Range("C20").WrapText = $true
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C20:D28").PasteSpecial(-4122) | Out-Null

# --- Row heights for the new rows ---
$ws.Range("A20:A28").RowHeight = 16.4

# --- Column widths for the new descr columns ---
$ws.Columns.Item(3).ColumnWidth = 48.00666666666667
$ws.Columns.Item(4).ColumnWidth = 45.696666666666665

# --- Update the view: scroll down and select the newly added block ---
$ws.Activate()
$excel.Goto($ws.Range("A10"), $false) | Out-Null
$ws.Range("C20:D28").Select() | Out-Null

Write-Output "Arabic language rows added"
